$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (autogluon): fill in previously empty cells
$ws.Range("B3").Value = "0.101 (0.054 ± 0.025)"
$ws.Range("C3").Value = "00:03:57 (00:04:38 ± 00:00:18)"
$ws.Range("D3").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("E3").Value = "[]"
# F3 holds a numeric-looking label ("61") that must stay a text value,
# matching how the sibling "best_seed" cells (F4="71", F6="7") are stored.
$f3 = $ws.Range("F3")
$f3.NumberFormat = "@"
$f3.Value = "61"
$f3.Style = "Normal"

# Row 4 (autokeras): fix mojibake "Â±" -> "±"
$ws.Range("B4").Value = "0.332 (0.253 ± 0.033)"
$ws.Range("C4").Value = "00:00:45 (00:00:55 ± 00:00:06)"
$ws.Range("D4").Value = "00:00:03 (00:00:03 ± 00:00:00)"

# Row 6 (autosklearn): fix mojibake "Â±" -> "±"
$ws.Range("B6").Value = "0.679 (0.565 ± 0.068)"
$ws.Range("C6").Value = "00:04:57 (00:05:01 ± 00:00:03)"
$ws.Range("D6").Value = "00:00:00 (00:00:06 ± 00:00:04)"
